# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.935.70'
$ws.Range("E2").Value = '  +2.06%  '

$ws.Range("D3").Value = '1.653.06'
$ws.Range("E3").Value = '  +2.89%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.88'
$ws.Range("E5").Value = '  +1.20%  '

$ws.Range("E6").Value = '  +2.23%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  +2.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0615'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.20'
$ws.Range("E10").Value = '  +4.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0874'
$ws.Range("E11").Value = '  +2.10%  '

$ws.Range("D12").Value = '1.888.26'
$ws.Range("E12").Value = '  +3.00%  '

$ws.Range("D13").Value = '1.655.60'
$ws.Range("E13").Value = '  +3.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("E14").Value = '  +1.81%  '

$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.07'
$ws.Range("E16").Value = '  +2.60%  '

$ws.Range("D17").Value = '26.937.97'
$ws.Range("E17").Value = '  +2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.97'
$ws.Range("E18").Value = '  +1.51%  '

$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.73'
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").Value = '  +3.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.29'
$ws.Range("E23").Value = '  +3.51%  '

$ws.Range("E24").Value = '  +3.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.29'
$ws.Range("E25").Value = '  -1.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.09'
$ws.Range("E26").Value = '  +1.79%  '

$ws.Range("E27").Value = '  +0.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.79'
$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("E30").Value = '  +0.34%  '

$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("D32").Value = '1.537.33'
$ws.Range("E32").Value = '  +3.20%  '

$ws.Range("E33").Value = '  +2.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").Value = '  +4.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +8.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.584'
$ws.Range("E37").Value = '  +3.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.890'
$ws.Range("E38").Value = '  +8.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0169'
$ws.Range("E39").Value = '  +2.77%  '

$ws.Range("E40").Value = '  +3.09%  '

$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("E42").Value = '  +2.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.70'
$ws.Range("E43").Value = '  +7.79%  '

$ws.Range("D44").Value = '1.794.73'
$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.775'
$ws.Range("E45").Value = '  +1.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.923'
$ws.Range("E46").Value = '  -1.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.28'
$ws.Range("E47").Value = '  +1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.51'
$ws.Range("E48").Value = '  +1.41%  '

$ws.Range("E49").Value = '  -1.65%  '

$ws.Range("E50").Value = '  +2.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0505'
$ws.Range("E51").Value = '  +0.89%  '
